# "Placering af svar boxe" - reposition the answer/question boxes on slide 1.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "A" answer-box (shape "Rectangle 4", id=5): nudge down (Top 140.7079pt -> 144.7079pt,
# i.e. y="1786990" EMU -> y="1837790" EMU).
$answerA = $s.Shapes.Item("Rectangle 4")
$answerA.Top = 144.70788

# Question box (shape "Rectangle 9", id=10): nudge up (Top 21.6pt -> 13.6pt,
# i.e. y="274320" EMU -> y="172720" EMU).
$question = $s.Shapes.Item("Rectangle 9")
$question.Top = 13.6
